$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ratingAmount/ratingValue for row 242
$ws.Range("D242").Value = 2
$ws.Range("E242").Value = 3

# Update timestamp column (O) for all data rows (2..398) to new scrape time
$newTimestamp = "2023-01-06 12:56:11"
for ($r = 2; $r -le 398; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
